$d = $word.ActiveDocument

# The document contains one top-level 3x3 table. Its second row, second
# column cell holds a nested 1x1 bordered table followed by an empty
# paragraph. We need to remove that nested table while keeping the
# trailing empty paragraph in the cell.
$outer = $d.Tables.Item(1)
$cell = $outer.Cell(2, 2)

# Replacing the cell's range contents with empty XML clears out the
# nested table (and any runs) while Word still preserves the cell's own
# paragraph mark, leaving a single empty paragraph behind - matching the
# target document exactly.
$cell.Range.InsertXML("")
